# Commit: "change all units to SI/pint-compatible"
# This workbook has a "Unit" column (E) containing "EJ/y" for the two
# data rows. Convert that unit string to the SI/pint-compatible "EJ".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the "Unit" header column so we don't hard-code positions unnecessarily.
$headerRow = 1
$unitCol = $null
for ($c = 1; $c -le $ws.UsedRange.Columns.Count; $c++) {
    $val = $ws.Cells.Item($headerRow, $c).Value()
    if ($val -eq "Unit") {
        $unitCol = $c
    }
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $unitCol)
    if ($cell.Value() -eq "EJ/y") {
        $cell.Value = "EJ"
    }
}

# Reflect the final selection recorded in the saved file: the merged
# notes cell at the bottom of the sheet (A4:H4).
$ws.Range("A4:H4").Select()

$wb.Save()
